# Apply the 2022-06-07 data update for Fonds de solidarite (VOLET1)
# Each entry updates the "nombre_aides" (column C) and "montant_total" (column E) values
# for a given row, leaving "nombre_entreprises" (column D) unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 10;  C = 345542;  E = 1817724469 },
    @{ Row = 35;  C = 6974;    E = 32982760 },
    @{ Row = 36;  C = 211205;  E = 404243995 },
    @{ Row = 67;  C = 27103;   E = 168710610 },
    @{ Row = 96;  C = 29544;   E = 56437349 },
    @{ Row = 100; C = 9345;    E = 23837942 },
    @{ Row = 103; C = 468;     E = 996792 },
    @{ Row = 121; C = 1306296; E = 2275168996 },
    @{ Row = 129; C = 633650;  E = 3432321719 },
    @{ Row = 132; C = 585914;  E = 3470178008 },
    @{ Row = 144; C = 25083;   E = 92541885 },
    @{ Row = 178; C = 515886;  E = 891200853 },
    @{ Row = 186; C = 236832;  E = 1189961448 },
    @{ Row = 237; C = 283321;  E = 1438426190 },
    @{ Row = 240; C = 205918;  E = 1069509020 }
)

foreach ($update in $updates) {
    $ws.Cells.Item($update.Row, 3).Value = $update.C
    $ws.Cells.Item($update.Row, 5).Value = $update.E
}
